# Lecture partielle de l'EDT M1 MIAGE.
#
# The timetable is shifted forward by exactly 3 calendar years: every date
# in column A moves from 2023 to 2026 (a +1096 day offset, which spans the
# 2024-02-29 leap day), and the French weekday label stored next to each
# date in column B is updated to match the weekday the shifted date now
# falls on. Two mistyped time slots (in the "CPO" / "BD" rows for the week
# of 2023-03-26) are also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;   Day = "lundi" },
    @{ Row = 4;   Day = "mercredi" },
    @{ Row = 6;   Day = "jeudi" },
    @{ Row = 10;  Day = "vendredi" },
    @{ Row = 13;  Day = "lundi" },
    @{ Row = 17;  Day = "mardi" },
    @{ Row = 19;  Day = "mercredi" },
    @{ Row = 22;  Day = "jeudi" },
    @{ Row = 26;  Day = "vendredi" },
    @{ Row = 32;  Day = "lundi" },
    @{ Row = 38;  Day = "mardi" },
    @{ Row = 40;  Day = "mercredi" },
    @{ Row = 42;  Day = "vendredi" },
    @{ Row = 48;  Day = "lundi" },
    @{ Row = 51;  Day = "mardi" },
    @{ Row = 55;  Day = "mercredi" },
    @{ Row = 57;  Day = "jeudi" },
    @{ Row = 62;  Day = "vendredi" },
    @{ Row = 66;  Day = "lundi" },
    @{ Row = 69;  Day = "mardi" },
    @{ Row = 73;  Day = "jeudi" },
    @{ Row = 78;  Day = "vendredi" },
    @{ Row = 81;  Day = "lundi" },
    @{ Row = 84;  Day = "mardi" },
    @{ Row = 86;  Day = "vendredi" },
    @{ Row = 90;  Day = "lundi" },
    @{ Row = 94;  Day = "lundi" },
    @{ Row = 97;  Day = "mardi" },
    @{ Row = 100; Day = "jeudi" },
    @{ Row = 102; Day = "vendredi" },
    @{ Row = 105; Day = "mardi" },
    @{ Row = 109; Day = "mercredi" },
    @{ Row = 112; Day = "jeudi" },
    @{ Row = 115; Day = "vendredi" },
    @{ Row = 120; Day = "lundi" },
    @{ Row = 122; Day = "mardi" },
    @{ Row = 125; Day = "mercredi" },
    @{ Row = 127; Day = "jeudi" },
    @{ Row = 130; Day = "vendredi" },
    @{ Row = 136; Day = "mardi" },
    @{ Row = 139; Day = "mercredi" },
    @{ Row = 141; Day = "jeudi" },
    @{ Row = 144; Day = "vendredi" },
    @{ Row = 148; Day = "lundi" },
    @{ Row = 151; Day = "mardi" },
    @{ Row = 153; Day = "mercredi" },
    @{ Row = 155; Day = "jeudi" },
    @{ Row = 158; Day = "vendredi" },
    @{ Row = 161; Day = "jeudi" },
    @{ Row = 164; Day = "mardi" },
    @{ Row = 166; Day = "jeudi" },
    @{ Row = 169; Day = "jeudi" },
    @{ Row = 172; Day = "mardi" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1096
    $ws.Cells.Item($r, 2).Value = $item.Day
}

$ws.Cells.Item(142, 4).Value = "10:0"
$ws.Cells.Item(143, 4).Value = "13:30"
